$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# The "choices" sheet lists Assistant codes in rows 2-16 (list_name="Assistant"),
# followed by a blank separator row, then other choice lists. Add a new
# assistant "Suzete" (code SUZ) as the last row of the Assistant block,
# i.e. insert a new row right before the current row 17 (the blank
# separator row under the Assistant list), shifting everything below down.
$ws.Range("A17").EntireRow.Insert()

$ws.Cells.Item(17, 1).Value = "Assistant"
$ws.Cells.Item(17, 2).Value = "SUZ"
$ws.Cells.Item(17, 3).Value = "Suzete"
$ws.Cells.Item(17, 4).Value = "Suzete"

# Newly materialized cells in columns C/D pick up the column's default
# style explicitly; clear back to Normal so the row matches the
# unstyled look of the other rows in this list (same as row 16 above it).
$ws.Range("C17:D17").Style = "Normal"

# Reflect the editing session ending with the "choices" tab active and
# the cursor on E14 (matches the user having just edited this sheet).
$ws.Activate()
$ws.Range("E14").Select()
